$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-21, columns A-F (data was reshuffled / reordered)
$data = @{
    2  = @(1001, 18, 30, 75, 60, 72)
    3  = @(901, 16, 15, 45, 60, 60)
    4  = @(801, 3, 67, 65, 52, 45)
    5  = @(1201, 2, 10, 10, 10, 10)
    6  = @(101, 9, 30, 15, 60, 15)
    7  = @(301, 6, 45, 30, 60, 45)
    8  = @(401, 9, 48, 67, 75, 45)
    9  = @(601, 9, 60, 67, 60, 42)
    10 = @(1202, 2, 10, 10, 10, 10)
    11 = @(902, 1, 0, 0, 0, 0)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(201, 9, 30, 15, 45, 30)
    15 = @(1203, 3, 15, 15, 15, 15)
    16 = @(1101, 0, 15, 30, 30, 0)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(502, 0, 4, 0, 0, 0)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(2, 0, 2, 2, 2, 2)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
